# Update gh-pages to output generated at 456a3b4
# Applies refreshed "remaining tickets" (F column) counts and a couple of
# sold-out -> not-for-sale (G column) status changes across the four
# worksheets of the workbook.

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibition) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 295
$ws.Range("G2").Value = "不可售"
$ws.Range("F4").Value = 257
$ws.Range("F5").Value = 726
$ws.Range("F6").Value = 2172
$ws.Range("F7").Value = 223
$ws.Range("F8").Value = 665
$ws.Range("F9").Value = 31
$ws.Range("F10").Value = 186
$ws.Range("F12").Value = 661
$ws.Range("F13").Value = 51
$ws.Range("F14").Value = 91
$ws.Range("F15").Value = 1298
$ws.Range("F17").Value = 44
$ws.Range("F18").Value = 188

# ---- 演出 (Performance) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 26
$ws.Range("F6").Value = 11
$ws.Range("F11").Value = 33

# ---- 本地生活 (Local Life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6326
$ws.Range("F3").Value = 787
$ws.Range("F4").Value = 1981
$ws.Range("F5").Value = 211

# ---- 全部类型 (All Types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6326
$ws.Range("F3").Value = 787
$ws.Range("F4").Value = 1981
$ws.Range("F5").Value = 295
$ws.Range("G5").Value = "不可售"
$ws.Range("F6").Value = 211
$ws.Range("F9").Value = 26
$ws.Range("F12").Value = 257
$ws.Range("F13").Value = 729
$ws.Range("F14").Value = 11
$ws.Range("F17").Value = 2172
$ws.Range("F19").Value = 223
$ws.Range("F21").Value = 33
$ws.Range("F22").Value = 665
$ws.Range("F23").Value = 31
$ws.Range("F24").Value = 186
$ws.Range("F27").Value = 661
$ws.Range("F28").Value = 51
$ws.Range("F29").Value = 91
$ws.Range("F31").Value = 1298
$ws.Range("F35").Value = 44
$ws.Range("F36").Value = 188
